$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) column for all existing data rows (2..133)
#    from 45192 to 45202.
for ($r = 2; $r -le 133; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2) Row 133 picks up an explicit custom row height (15pt) in the edit.
$ws.Rows.Item(133).RowHeight = 15

# 3) Append the new record as row 134.
$ws.Cells.Item(134, 1).Value = "A 46389-2023"
$ws.Cells.Item(134, 2).Value = 45197
$ws.Cells.Item(134, 3).Value = 45202
$ws.Cells.Item(134, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(134, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(134, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item(134, 5).Value = "STORFORS"
$ws.Cells.Item(134, 7).Value = 3.1
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = 0
$ws.Cells.Item(134, 14).Value = 0
$ws.Cells.Item(134, 15).Value = 0
$ws.Cells.Item(134, 16).Value = 0
$ws.Cells.Item(134, 17).Value = 0

# R column keeps the same "wrap text" style used throughout the sheet,
# with an empty (inline) string value.
$ws.Cells.Item(134, 18).Value = ""
$ws.Cells.Item(134, 18).WrapText = $true
